$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Hyperlinks.Add($ws.Range("B5"), "http://china.rs-online.com/web/p/general-purpose-adcs/7094550/")

$ws.Range("B5").Select()
